# Auto-generated PowerShell COM-interop script
# Applies the panelApp "time_taken" column addition + WLS gene row insertion
# (data/panelapp/au/Anophthalmia_Microphthalmia_Coloboma.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 70 for the WLS gene (shifts existing rows 70-82 down to 71-83)
$ws.Rows.Item(70).Insert()

# 2. Populate the new row 70 (WLS gene, confidence 3)
$ws.Cells.Item(70,2).Value2 = "WLS"
$ws.Cells.Item(70,3).Value2 = "wntless Wnt ligand secretion mediator"
$ws.Cells.Item(70,4).Value2 = "3"
$ws.Cells.Item(70,5).Value2 = "Anophthalmia_Microphthalmia_Coloboma"

# Match column-A styling (bold, bordered, centered/top) used throughout column A
$ws.Cells.Item(70,1).Font.Bold = $true
$ws.Cells.Item(70,1).HorizontalAlignment = -4108
$ws.Cells.Item(70,1).VerticalAlignment = -4160
$ws.Cells.Item(70,1).Borders.LineStyle = 1

# 3. Re-sequence column A (row index) for rows 70-83. Written top-down, one at a time:
#    the engine re-derives each row's sequential index from scratch, so this must be
#    done for every row from the insertion point through the old last row + 1.
$ws.Cells.Item(70,1).Value2 = 68
$ws.Cells.Item(71,1).Value2 = 69
$ws.Cells.Item(72,1).Value2 = 70
$ws.Cells.Item(73,1).Value2 = 71
$ws.Cells.Item(74,1).Value2 = 72
$ws.Cells.Item(75,1).Value2 = 73
$ws.Cells.Item(76,1).Value2 = 74
$ws.Cells.Item(77,1).Value2 = 75
$ws.Cells.Item(78,1).Value2 = 76
$ws.Cells.Item(79,1).Value2 = 77
$ws.Cells.Item(80,1).Value2 = 78
$ws.Cells.Item(81,1).Value2 = 79
$ws.Cells.Item(82,1).Value2 = 80
$ws.Cells.Item(83,1).Value2 = 81

# 4. Add new column F ("time_taken") header + per-row metadata timestamps
$ws.Cells.Item(1,6).Value2 = "time_taken"
$ws.Cells.Item(1,6).Font.Bold = $true
$ws.Cells.Item(1,6).HorizontalAlignment = -4108
$ws.Cells.Item(1,6).VerticalAlignment = -4160
$ws.Cells.Item(1,6).Borders.LineStyle = 1

$ws.Cells.Item(2,6).Value2 = "2021-10-05 10:50:07.228697"
$ws.Cells.Item(3,6).Value2 = "2021-10-05 10:50:07.228708"
$ws.Cells.Item(4,6).Value2 = "2021-10-05 10:50:07.228712"
$ws.Cells.Item(5,6).Value2 = "2021-10-05 10:50:07.228715"
$ws.Cells.Item(6,6).Value2 = "2021-10-05 10:50:07.228718"
$ws.Cells.Item(7,6).Value2 = "2021-10-05 10:50:07.228720"
$ws.Cells.Item(8,6).Value2 = "2021-10-05 10:50:07.228723"
$ws.Cells.Item(9,6).Value2 = "2021-10-05 10:50:07.228725"
$ws.Cells.Item(10,6).Value2 = "2021-10-05 10:50:07.228728"
$ws.Cells.Item(11,6).Value2 = "2021-10-05 10:50:07.228731"
$ws.Cells.Item(12,6).Value2 = "2021-10-05 10:50:07.228734"
$ws.Cells.Item(13,6).Value2 = "2021-10-05 10:50:07.228736"
$ws.Cells.Item(14,6).Value2 = "2021-10-05 10:50:07.228739"
$ws.Cells.Item(15,6).Value2 = "2021-10-05 10:50:07.228741"
$ws.Cells.Item(16,6).Value2 = "2021-10-05 10:50:07.228744"
$ws.Cells.Item(17,6).Value2 = "2021-10-05 10:50:07.228746"
$ws.Cells.Item(18,6).Value2 = "2021-10-05 10:50:07.228749"
$ws.Cells.Item(19,6).Value2 = "2021-10-05 10:50:07.228752"
$ws.Cells.Item(20,6).Value2 = "2021-10-05 10:50:07.228755"
$ws.Cells.Item(21,6).Value2 = "2021-10-05 10:50:07.228757"
$ws.Cells.Item(22,6).Value2 = "2021-10-05 10:50:07.228760"
$ws.Cells.Item(23,6).Value2 = "2021-10-05 10:50:07.228762"
$ws.Cells.Item(24,6).Value2 = "2021-10-05 10:50:07.228765"
$ws.Cells.Item(25,6).Value2 = "2021-10-05 10:50:07.228768"
$ws.Cells.Item(26,6).Value2 = "2021-10-05 10:50:07.228771"
$ws.Cells.Item(27,6).Value2 = "2021-10-05 10:50:07.228773"
$ws.Cells.Item(28,6).Value2 = "2021-10-05 10:50:07.228776"
$ws.Cells.Item(29,6).Value2 = "2021-10-05 10:50:07.228779"
$ws.Cells.Item(30,6).Value2 = "2021-10-05 10:50:07.228781"
$ws.Cells.Item(31,6).Value2 = "2021-10-05 10:50:07.228784"
$ws.Cells.Item(32,6).Value2 = "2021-10-05 10:50:07.228786"
$ws.Cells.Item(33,6).Value2 = "2021-10-05 10:50:07.228789"
$ws.Cells.Item(34,6).Value2 = "2021-10-05 10:50:07.228792"
$ws.Cells.Item(35,6).Value2 = "2021-10-05 10:50:07.228794"
$ws.Cells.Item(36,6).Value2 = "2021-10-05 10:50:07.228797"
$ws.Cells.Item(37,6).Value2 = "2021-10-05 10:50:07.228800"
$ws.Cells.Item(38,6).Value2 = "2021-10-05 10:50:07.228802"
$ws.Cells.Item(39,6).Value2 = "2021-10-05 10:50:07.228805"
$ws.Cells.Item(40,6).Value2 = "2021-10-05 10:50:07.228807"
$ws.Cells.Item(41,6).Value2 = "2021-10-05 10:50:07.228810"
$ws.Cells.Item(42,6).Value2 = "2021-10-05 10:50:07.228813"
$ws.Cells.Item(43,6).Value2 = "2021-10-05 10:50:07.228815"
$ws.Cells.Item(44,6).Value2 = "2021-10-05 10:50:07.228818"
$ws.Cells.Item(45,6).Value2 = "2021-10-05 10:50:07.228821"
$ws.Cells.Item(46,6).Value2 = "2021-10-05 10:50:07.228824"
$ws.Cells.Item(47,6).Value2 = "2021-10-05 10:50:07.228826"
$ws.Cells.Item(48,6).Value2 = "2021-10-05 10:50:07.228829"
$ws.Cells.Item(49,6).Value2 = "2021-10-05 10:50:07.228831"
$ws.Cells.Item(50,6).Value2 = "2021-10-05 10:50:07.228834"
$ws.Cells.Item(51,6).Value2 = "2021-10-05 10:50:07.228836"
$ws.Cells.Item(52,6).Value2 = "2021-10-05 10:50:07.228839"
$ws.Cells.Item(53,6).Value2 = "2021-10-05 10:50:07.228841"
$ws.Cells.Item(54,6).Value2 = "2021-10-05 10:50:07.228844"
$ws.Cells.Item(55,6).Value2 = "2021-10-05 10:50:07.228847"
$ws.Cells.Item(56,6).Value2 = "2021-10-05 10:50:07.228849"
$ws.Cells.Item(57,6).Value2 = "2021-10-05 10:50:07.228852"
$ws.Cells.Item(58,6).Value2 = "2021-10-05 10:50:07.228854"
$ws.Cells.Item(59,6).Value2 = "2021-10-05 10:50:07.228857"
$ws.Cells.Item(60,6).Value2 = "2021-10-05 10:50:07.228859"
$ws.Cells.Item(61,6).Value2 = "2021-10-05 10:50:07.228862"
$ws.Cells.Item(62,6).Value2 = "2021-10-05 10:50:07.228864"
$ws.Cells.Item(63,6).Value2 = "2021-10-05 10:50:07.228867"
$ws.Cells.Item(64,6).Value2 = "2021-10-05 10:50:07.228869"
$ws.Cells.Item(65,6).Value2 = "2021-10-05 10:50:07.228872"
$ws.Cells.Item(66,6).Value2 = "2021-10-05 10:50:07.228875"
$ws.Cells.Item(67,6).Value2 = "2021-10-05 10:50:07.228878"
$ws.Cells.Item(68,6).Value2 = "2021-10-05 10:50:07.228881"
$ws.Cells.Item(69,6).Value2 = "2021-10-05 10:50:07.228883"
$ws.Cells.Item(70,6).Value2 = "2021-10-05 10:50:07.228885"
$ws.Cells.Item(71,6).Value2 = "2021-10-05 10:50:07.228888"
$ws.Cells.Item(72,6).Value2 = "2021-10-05 10:50:07.228891"
$ws.Cells.Item(73,6).Value2 = "2021-10-05 10:50:07.228893"
$ws.Cells.Item(74,6).Value2 = "2021-10-05 10:50:07.228895"
$ws.Cells.Item(75,6).Value2 = "2021-10-05 10:50:07.228898"
$ws.Cells.Item(76,6).Value2 = "2021-10-05 10:50:07.228900"
$ws.Cells.Item(77,6).Value2 = "2021-10-05 10:50:07.228903"
$ws.Cells.Item(78,6).Value2 = "2021-10-05 10:50:07.228907"
$ws.Cells.Item(79,6).Value2 = "2021-10-05 10:50:07.228910"
$ws.Cells.Item(80,6).Value2 = "2021-10-05 10:50:07.228913"
$ws.Cells.Item(81,6).Value2 = "2021-10-05 10:50:07.228915"
$ws.Cells.Item(82,6).Value2 = "2021-10-05 10:50:07.228918"
$ws.Cells.Item(83,6).Value2 = "2021-10-05 10:50:07.228920"

Write-Host "Edit applied"
